# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master / applied design
#                             ("Integral" / "Red Violet" colours)
#   ppt/theme/theme2.xml  -> bound to the notes master
#                             ("Office Theme" / "Office" colours)
#
# The authored change swaps the two theme bodies, so the design that is
# actually applied to the deck (theme1.xml) becomes the stock "Office
# Theme" colour palette, and the notes-master theme becomes the palette
# that used to be applied (Integral/Red Violet).
#
# Drive this the same way a user would from the Design tab: push the
# target theme's 12 colour-scheme entries onto the presentation's applied
# theme via the SlideMaster's Theme object.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$scheme = $m.Theme.ThemeColorScheme

# Index order for ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
# Values are OLE COLORREF ints (0x00BBGGRR) for the "Office Theme" palette.
$scheme.Item(1).RGB  = 0        # dk1      000000
$scheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388  # dk2      44546A
$scheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407    # accent4  FFC000
$scheme.Item(9).RGB  = 12874308 # accent5  4472C4
$scheme.Item(10).RGB = 4697456  # accent6  70AD47
$scheme.Item(11).RGB = 12673797 # hlink    0563C1
$scheme.Item(12).RGB = 7491477  # folHlink 954F72
